$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. ListConfig: add a 4th column "RunMode" with Y/N list values,
#    used as the source list for the new TestCases RunMode dropdown.
# ------------------------------------------------------------------
$listConfig = $wb.Worksheets.Item("ListConfig")
$listConfig.Range("D1").Value = "Y"
$listConfig.Range("D2").Value = "N"
$listConfig.Range("D1:D2").HorizontalAlignment = -4108
$listConfig.Range("D1:D2").VerticalAlignment = -4108
$listConfig.Range("D1").Select()

# ------------------------------------------------------------------
# 2. Insert a new "TestCases" worksheet between TestData and
#    ListConfig, holding the generic test-case names and whether
#    each one should run.
# ------------------------------------------------------------------
$testData = $wb.Worksheets.Item("TestData")
$testCases = $wb.Worksheets.Add($null, $testData)
$testCases.Name = "TestCases"

$testCases.Range("A1").Value = "TestNames"
$testCases.Range("B1").Value = "RunMode"
$testCases.Range("A2").Value = "LoginAsBankManager"
$testCases.Range("B2").Value = "Y"
$testCases.Range("A3").Value = "AddCustomers"
$testCases.Range("B3").Value = "Y"
$testCases.Range("A4").Value = "OpenAccount"
$testCases.Range("B4").Value = "N"

# Header formatting -- same bold / yellow-filled / bordered header
# style used on the other sheets of the workbook.
$header = $testCases.Range("A1:B1")
$header.Font.Bold = $true
$header.Interior.Color = 65535
$header.Borders.LineStyle = 1

# Body formatting: bordered cells, RunMode column centered.
$bodyA = $testCases.Range("A2:A4")
$bodyA.Borders.LineStyle = 1

$bodyB = $testCases.Range("B2:B4")
$bodyB.Borders.LineStyle = 1
$bodyB.HorizontalAlignment = -4108
$bodyB.VerticalAlignment = -4108

# Column widths to fit the test-case names / run-mode values.
$testCases.Columns.Item(1).ColumnWidth = 19.29
$testCases.Columns.Item(2).ColumnWidth = 10.29

# Data validation: RunMode must be Y/N, sourced from ListConfig!D1:D2
$testCases.Range("B2:B4").Validation.Add(3, 1, 1, "=ListConfig!`$D`$1:`$D`$2")

$testCases.Range("A1").Select()
